$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 76, pushing existing rows 76-114 down to 77-115
$ws.Rows("76").Insert()

# Populate the new row 76 with the new weekly record
$ws.Range("A76").Value = 10
$ws.Range("B76").Value = "Vega Modelo de Temuco"
$ws.Range("C76").Value = "La Araucanía"
$ws.Range("D76").Value = (Get-Date -Year 2023 -Month 11 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E76").Value = 9
$ws.Range("F76").Value = 300000000
$ws.Range("G76").Value = "Espárragos"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 300
$ws.Range("K76").Value = 1800
$ws.Range("L76").Value = 1800
$ws.Range("M76").Value = 1800
$ws.Range("N76").Value = "`$/kilo"
$ws.Range("O76").Value = "Región del Maule"
$ws.Range("P76").Value = 1800
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"
